$d = $word.ActiveDocument
$d.Content.Find.Execute("فارسی", $true, $false, $false, $false, $false, $true, 1, $false, "پارسی", 2)
